# Femacal de La Calera - Alcachofa: insert two new weekly price rows
# (Provincia de Quillota, Española, 09-11-2021) above the existing
# "Primera"/Argentina(o) row dated 13-08-2021, shifting the rest of the
# table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 237 (pushes old rows 237-244 down to 239-246)
$ws.Rows.Item(237).Insert()
$ws.Rows.Item(237).Insert()

# New row 237
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44509
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112013
$ws.Cells.Item(237, 7).Value = "Alcachofa"
$ws.Cells.Item(237, 8).Value = "Española"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 12300
$ws.Cells.Item(237, 11).Value = 290
$ws.Cells.Item(237, 12).Value = 300
$ws.Cells.Item(237, 13).Value = 296
$ws.Cells.Item(237, 14).Value = "$/unidad"
$ws.Cells.Item(237, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(237, 16).Value = 296
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"

# New row 238
$ws.Cells.Item(238, 1).Value = 3
$ws.Cells.Item(238, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(238, 3).Value = "Coquimbo"
$ws.Cells.Item(238, 4).Value = 44509
$ws.Cells.Item(238, 5).Value = 5
$ws.Cells.Item(238, 6).Value = 100112013
$ws.Cells.Item(238, 7).Value = "Alcachofa"
$ws.Cells.Item(238, 8).Value = "Española"
$ws.Cells.Item(238, 9).Value = "Segunda"
$ws.Cells.Item(238, 10).Value = 6500
$ws.Cells.Item(238, 11).Value = 200
$ws.Cells.Item(238, 12).Value = 200
$ws.Cells.Item(238, 13).Value = 200
$ws.Cells.Item(238, 14).Value = "$/unidad"
$ws.Cells.Item(238, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(238, 16).Value = 200
$ws.Cells.Item(238, 17).Value = 1
$ws.Cells.Item(238, 18).Value = "Hortaliza"
